# Generate Report for Handback
# Updates the zh-cn and de-de localization-status sheets: the pending
# "Ready for handoff" row for 5e8f358b-....md now reports a failed
# handback transform, with the Error Detail column populated and widened.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status column (C) for the 5e8f358b-...md row (row 3) switches from
# "Ready for handoff" to "Handback transform failed" on both language
# sheets; the Overview sheet mirrors the same status text in its
# zh-cn/de-de columns (E3/F3) for that row.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Error Detail column (P) for that same row now carries the mismatch
# message describing the handback/handoff file-name discrepancy.
$zhcn.Range("P3").Value = "Handback file name: u5nk2zz2.z1a is different with handoff file name: 5e8f358b-fe43-4309-915b-7fffb7172ca6.d56ad895900b6159721de8a81171dbc10b6cc031.zh-cn."
$dede.Range("P3").Value = "Handback file name: u5nk2zz2.z1a is different with handoff file name: 5e8f358b-fe43-4309-915b-7fffb7172ca6.d56ad895900b6159721de8a81171dbc10b6cc031.de-de."

# Widen the Error Detail column (16 / P) to fit the new message text.
# ColumnWidth is stored in the package as (value + 5/6), so back off by
# 5/6 to land on an exact "40" in the saved width attribute.
$zhcn.Columns.Item(16).ColumnWidth = 40 - 5/6
$dede.Columns.Item(16).ColumnWidth = 40 - 5/6
